$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.031031
$ws.Range("H2").Value = 12.093093
$ws.Range("I2").Value = 0.380357182622003
$ws.Range("J2").Value = 0.380357182622003
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 10.980535162385
$ws.Range("R2").Value = 98.82481646146501
$ws.Range("S2").Value = 0.01763145900423589
$ws.Range("T2").Value = 0.01763145900423589
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.031031
$ws.Range("H3").Value = 12.093093
$ws.Range("I3").Value = 0.380357182622003
$ws.Range("J3").Value = 0.380357182622003
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 163.743021456884
$ws.Range("R3").Value = 1473.687193111956
$ws.Range("S3").Value = 0.2629223737597594
$ws.Range("T3").Value = 0.2629223737597594
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.031031
$ws.Range("H4").Value = 12.093093
$ws.Range("I4").Value = 0.380357182622003
$ws.Range("J4").Value = 0.380357182622003
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 62.155615832835
$ws.Range("R4").Value = 559.400542495515
$ws.Range("S4").Value = 0.09980334985800773
$ws.Range("T4").Value = 0.09980334985800773
$ws.Range("I5").Value = 0.4810839099297969
$ws.Range("J5").Value = 0.4810839099297969
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 13.88841602155722
$ws.Range("R5").Value = 124.995744194015
$ws.Range("S5").Value = 0.02230064692627168
$ws.Range("T5").Value = 0.02230064692627168
$ws.Range("I6").Value = 0.4810839099297969
$ws.Range("J6").Value = 0.4810839099297969
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("S6").Value = 0.3325498488142692
$ws.Range("T6").Value = 0.3325498488142692
$ws.Range("I7").Value = 0.4810839099297969
$ws.Range("J7").Value = 0.4810839099297969
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 78.61575396795166
$ws.Range("R7").Value = 707.541785711565
$ws.Range("S7").Value = 0.126233414189256
$ws.Range("T7").Value = 0.126233414189256
$ws.Range("G8").Value = 1.468449333333333
$ws.Range("H8").Value = 4.405348
$ws.Range("I8").Value = 0.1385589074482
$ws.Range("J8").Value = 0.1385589074482
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 4.000058431415555
$ws.Range("R8").Value = 36.00052588274
$ws.Range("S8").Value = 0.006422898811858355
$ws.Range("T8").Value = 0.006422898811858355
$ws.Range("G9").Value = 1.468449333333333
$ws.Range("H9").Value = 4.405348
$ws.Range("I9").Value = 0.1385589074482
$ws.Range("J9").Value = 0.1385589074482
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 59.64933802204622
$ws.Range("R9").Value = 536.8440421984161
$ws.Range("S9").Value = 0.09577901645160659
$ws.Range("T9").Value = 0.09577901645160658
$ws.Range("G10").Value = 1.468449333333333
$ws.Range("H10").Value = 4.405348
$ws.Range("I10").Value = 0.1385589074482
$ws.Range("J10").Value = 0.1385589074482
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 22.64243877872667
$ws.Range("R10").Value = 203.78194900854
$ws.Range("S10").Value = 0.03635699218473509
$ws.Range("T10").Value = 0.03635699218473509
